$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil2")
$ws.Range("C2").Value = "http://data.sparna.fr/vocabularies/days/"
$ws.Activate()
$ws.Range("C9").Select()
